$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# are temporarily switched to Text format, written, then restored to the
# default "Normal" style so the final style index matches the original (0).

$ws.Range("D2").Value = "23.396.13"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.630.19"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9996"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "302.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3769"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.94"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3629"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = "  -2.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9995"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.469"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.325"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("D17").Value = "1.619.97"
$ws.Range("E17").Value = "  -0.79%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.538"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("D22").Style = "Normal"
$ws.Range("E23").Value = "  -2.12%  "
$ws.Range("D24").Value = "23.391.71"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.510"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.06%  "
$ws.Range("E26").Value = "  -0.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.274"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.37%  "
$ws.Range("D31").Value = "1.803.15"
$ws.Range("E31").Value = "  -0.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.614"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.138"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.059"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02764"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2489"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08756"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.07133"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.966"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6986"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.330"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "15.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.30%  "
$ws.Range("E44").Value = "  -3.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6458"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9989"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.274"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.959"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07972"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "126.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("E51").Value = "  -1.08%  "
